# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price (column D) and 1h volume change (column E) are stored as plain text,
# so values that look fully numeric (e.g. "523.12") are written with a
# leading apostrophe to keep Excel from auto-coercing them to numbers and
# losing the exact original formatting (trailing zeros, etc.). Values that
# already contain a second "." (e.g. "58.548.84") or other non-numeric
# characters are left as-is since Excel keeps those as text automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.548.84"
$ws.Range("E2").Value = "  +1.04%  "

$ws.Range("D3").Value = "2.523.86"
$ws.Range("E3").Value = "  +2.55%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'523.12"
$ws.Range("E5").Value = "  +0.96%  "

$ws.Range("D6").Value = "'133.29"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("D8").Value = "'0.560"
$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("D9").Value = "2.522.52"
$ws.Range("E9").Value = "  +2.25%  "

$ws.Range("D10").Value = "'0.0978"
$ws.Range("E10").Value = "  +0.44%  "

$ws.Range("E11").Value = "  -1.33%  "

$ws.Range("D12").Value = "'5.18"
$ws.Range("E12").Value = "  -2.11%  "

$ws.Range("E13").Value = "  -0.70%  "

$ws.Range("D14").Value = "2.971.81"
$ws.Range("E14").Value = "  +2.53%  "

$ws.Range("D15").Value = "58.519.12"
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("D16").Value = "'22.15"
$ws.Range("E16").Value = "  +1.12%  "

$ws.Range("E17").Value = "  +0.87%  "

$ws.Range("D18").Value = "2.520.59"
$ws.Range("E18").Value = "  +2.41%  "

$ws.Range("D19").Value = "'10.67"
$ws.Range("E19").Value = "  +0.96%  "

$ws.Range("D20").Value = "'322.05"
$ws.Range("E20").Value = "  +1.14%  "

$ws.Range("E21").Value = "  +0.96%  "

$ws.Range("D22").Value = "'6.17"
$ws.Range("E22").Value = "  +8.49%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").Value = "'64.77"
$ws.Range("E24").Value = "  +0.62%  "

$ws.Range("D25").Value = "'0.407"
$ws.Range("E25").Value = "  +0.22%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("E27").Value = "  +0.62%  "

$ws.Range("D28").Value = "'7.40"
$ws.Range("E28").Value = "  +1.41%  "

$ws.Range("D29").Value = "0.0₃0757"
$ws.Range("E29").Value = "  +2.38%  "

$ws.Range("D30").Value = "'1.73"
$ws.Range("E30").Value = "  +2.54%  "

$ws.Range("E31").Value = "  +3.66%  "

$ws.Range("D32").Value = "'167.87"
$ws.Range("E32").Value = "  -0.22%  "

$ws.Range("D33").Value = "'6.37"
$ws.Range("E33").Value = "  +2.66%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").Value = "'18.14"
$ws.Range("E36").Value = "  +0.95%  "

$ws.Range("E37").Value = "  -4.08%  "

$ws.Range("D38").Value = "'3.95"
$ws.Range("E38").Value = "  +0.40%  "

$ws.Range("E39").Value = "  +2.23%  "

$ws.Range("D40").Value = "'36.37"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").Value = "'0.775"
$ws.Range("E41").Value = "  -1.48%  "

$ws.Range("D42").Value = "'278.99"
$ws.Range("E42").Value = "  +3.13%  "

$ws.Range("D43").Value = "'3.50"
$ws.Range("E43").Value = "  +2.56%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'130.69"
$ws.Range("E44").Value = "  +5.89%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'5.01"
$ws.Range("E45").Value = "  +1.19%  "

$ws.Range("D46").Value = "'0.600"
$ws.Range("E46").Value = "  +2.07%  "

$ws.Range("D47").Value = "'0.0923"
$ws.Range("E47").Value = "  +1.86%  "

$ws.Range("D48").Value = "'0.0502"
$ws.Range("E48").Value = "  +3.95%  "

$ws.Range("D49").Value = "'17.80"
$ws.Range("E49").Value = "  +1.68%  "

$ws.Range("E50").Value = "  +1.51%  "

$ws.Range("D51").Value = "'16.94"
$ws.Range("E51").Value = "  +1.25%  "
